# Fix variables/query assignment on the "startup" sheet:
#  - B2 (CasesTab row) should hold the Cases query (with the stray
#    "Cohort" coalesce line removed)
#  - B3 (SamplesTab row) should hold the Sample query
#  - B4 (FilesTab row) should hold the Files query
# Afterwards the active selection moves from B4 up to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Irish Wolfhound''] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

$sampleQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN [''Irish Wolfhound''] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '''') AS `Sample ID`, 
        coalesce(c.case_id, '''') AS `Case ID`, 
        coalesce(demo.breed,'''') AS Breed , 
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(samp.sample_site, '''') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '''') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '''') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '''') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '''') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '''') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '''') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '''') AS `Sample Preservation`'

$filesQuery = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Irish Wolfhound''] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '''') AS `File Name`,
         coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`,
         coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'

$ws.Range("B2").Value = $casesQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $filesQuery

$ws.Range("B2").Select()
